$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3454506666666666
$ws.Range("H2").Value = 1.036352
$ws.Range("I2").Value = 0.1052716477644991
$ws.Range("J2").Value = 0.1052716477644991
$ws.Range("M2").Value = 6.875726333333334
$ws.Range("N2").Value = 20.627179
$ws.Range("O2").Value = 0.6245871044820662
$ws.Range("P2").Value = 0.6245871044820662
$ws.Range("Q2").Value = 2.375224245667555
$ws.Range("R2").Value = 21.377018211008
$ws.Range("S2").Value = 0.06575131366128449
$ws.Range("T2").Value = 0.06575131366128449

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3454506666666666
$ws.Range("H3").Value = 1.036352
$ws.Range("I3").Value = 0.1052716477644991
$ws.Range("J3").Value = 0.1052716477644991
$ws.Range("O3").Value = 0.06694469792011602
$ws.Range("P3").Value = 0.06694469792011602
$ws.Range("Q3").Value = 0.2545820566542222
$ws.Range("R3").Value = 2.291238509888
$ws.Range("S3").Value = 0.007047378659147252
$ws.Range("T3").Value = 0.007047378659147252

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3454506666666666
$ws.Range("H4").Value = 1.036352
$ws.Range("I4").Value = 0.1052716477644991
$ws.Range("J4").Value = 0.1052716477644991
$ws.Range("O4").Value = 0.3084681975978177
$ws.Range("P4").Value = 0.3084681975978177
$ws.Range("Q4").Value = 1.173064792234667
$ws.Range("R4").Value = 10.557583130112
$ws.Range("S4").Value = 0.03247295544406738
$ws.Range("T4").Value = 0.03247295544406739

# Row 5
$ws.Range("I5").Value = 0.3398937483175971
$ws.Range("J5").Value = 0.3398937483175971
$ws.Range("M5").Value = 6.875726333333334
$ws.Range("N5").Value = 20.627179
$ws.Range("O5").Value = 0.6245871044820662
$ws.Range("P5").Value = 0.6245871044820662
$ws.Range("Q5").Value = 7.668958253231001
$ws.Range("R5").Value = 69.020624279079
$ws.Range("S5").Value = 0.2122932520932441
$ws.Range("T5").Value = 0.2122932520932441

# Row 6
$ws.Range("I6").Value = 0.3398937483175971
$ws.Range("J6").Value = 0.3398937483175971
$ws.Range("O6").Value = 0.06694469792011602
$ws.Range("P6").Value = 0.06694469792011602
$ws.Range("S6").Value = 0.02275408430605748
$ws.Range("T6").Value = 0.02275408430605748

# Row 7
$ws.Range("I7").Value = 0.3398937483175971
$ws.Range("J7").Value = 0.3398937483175971
$ws.Range("O7").Value = 0.3084681975978177
$ws.Range("P7").Value = 0.3084681975978177
$ws.Range("Q7").Value = 3.787509720984001
$ws.Range("S7").Value = 0.1048464119182955
$ws.Range("T7").Value = 0.1048464119182955

# Row 8
$ws.Range("I8").Value = 0.5548346039179038
$ws.Range("J8").Value = 0.5548346039179038
$ws.Range("M8").Value = 6.875726333333334
$ws.Range("N8").Value = 20.627179
$ws.Range("O8").Value = 0.6245871044820662
$ws.Range("P8").Value = 0.6245871044820662
$ws.Range("Q8").Value = 12.51862805937367
$ws.Range("R8").Value = 112.667652534363
$ws.Range("S8").Value = 0.3465425387275377
$ws.Range("T8").Value = 0.3465425387275377

# Row 9
$ws.Range("I9").Value = 0.5548346039179038
$ws.Range("J9").Value = 0.5548346039179038
$ws.Range("O9").Value = 0.06694469792011602
$ws.Range("P9").Value = 0.06694469792011602
$ws.Range("S9").Value = 0.03714323495491129
$ws.Range("T9").Value = 0.03714323495491129

# Row 10
$ws.Range("I10").Value = 0.5548346039179038
$ws.Range("J10").Value = 0.5548346039179038
$ws.Range("O10").Value = 0.3084681975978177
$ws.Range("P10").Value = 0.3084681975978177
$ws.Range("R10").Value = 55.64378043583201
$ws.Range("S10").Value = 0.1711488302354549
$ws.Range("T10").Value = 0.1711488302354549
